$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: replace the worker's identification data ---
# (was CC / 1143397510 / KELLYS TATIANA ZAMBRANO MENDOZA / period 2507 / 2860 / 2144712)
$ws.Range("B16").Value = "CE"
$ws.Range("C16").Value = "2527776"
$ws.Range("D16").Value = "LINETH CAROLINA LOPEZ MEDINA"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 909000

# --- Remove the old row 17 (second mora period for the old worker) ---
# Row 18 (already LINETH CAROLINA LOPEZ MEDINA / 2527776 / CE) shifts up to become
# the new row 17, carrying its own formatting (the table's closing border).
$ws.Rows("17").Delete()

# The shifted-up row (now row 17) should reflect period 2508 instead of 2507.
$ws.Range("E17").Value = "2508"

# --- Update summary figures ---
# Cant. Trabajadores: 2 -> 1
$ws.Range("C13").Value = 1
# Valor Mora total: 62660 -> 113880
$ws.Range("E11").Value = 113880

# --- Column D autosize now that the (shorter) name replaced the old one ---
$ws.Columns("D").AutoFit()
